$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "14/03/2023"

$ws.Range("D3").Value = 270.4
$ws.Range("D4").Value = 361.4
$ws.Range("D5").Value = 377
$ws.Range("D6").Value = 331.5
$ws.Range("D7").Value = 234
$ws.Range("D8").Value = 284.7
$ws.Range("D9").Value = 258.7
$ws.Range("D10").Value = 301.6
$ws.Range("D11").Value = 280.8
$ws.Range("D12").Value = 308.1
$ws.Range("D13").Value = 230.1
$ws.Range("D14").Value = 150.8
$ws.Range("D18").Value = 9
$ws.Range("D19").Value = 34
$ws.Range("D20").Value = 44
$ws.Range("D21").Value = 33
$ws.Range("D22").Value = 41
$ws.Range("D23").Value = 30
$ws.Range("D24").Value = 36
$ws.Range("D25").Value = 25
$ws.Range("D26").Value = 37
$ws.Range("D27").Value = 25
$ws.Range("D28").Value = 39
$ws.Range("D29").Value = 16
$ws.Range("D30").Value = 22
$ws.Range("D34").Value = 46
$ws.Range("D35").Value = 66
$ws.Range("D36").Value = 82
$ws.Range("D37").Value = 165
$ws.Range("D38").Value = 218
$ws.Range("D39").Value = 213
$ws.Range("D40").Value = 179
$ws.Range("D41").Value = 148
$ws.Range("D42").Value = 166
$ws.Range("D43").Value = 145
$ws.Range("D44").Value = 153
$ws.Range("D45").Value = 162
$ws.Range("D46").Value = 155
$ws.Range("D47").Value = 129
$ws.Range("D48").Value = 78
$ws.Range("D49").Value = 48
$ws.Range("D50").Value = 23
$ws.Range("D52").Value = 75
$ws.Range("D53").Value = 135
$ws.Range("D54").Value = 168
$ws.Range("D55").Value = 158
$ws.Range("D56").Value = 131
$ws.Range("D57").Value = 110
$ws.Range("D58").Value = 127
$ws.Range("D59").Value = 114
$ws.Range("D60").Value = 130
$ws.Range("D61").Value = 114
$ws.Range("D62").Value = 128
$ws.Range("D63").Value = 93
$ws.Range("D68").Value = 27
$ws.Range("D69").Value = 46
$ws.Range("D70").Value = 47
$ws.Range("D71").Value = 45
$ws.Range("D72").Value = 49
$ws.Range("D73").Value = 44
$ws.Range("D74").Value = 41
$ws.Range("D75").Value = 47
$ws.Range("D76").Value = 44
$ws.Range("D77").Value = 53
$ws.Range("D78").Value = 41
$ws.Range("D79").Value = 33
$ws.Range("D80").Value = 18
$ws.Range("D81").Value = 11
$ws.Range("D82").Value = 5
$ws.Range("D87").Value = 1
$ws.Range("D88").Value = 2
$ws.Range("D89").Value = 3
$ws.Range("D90").Value = 16
$ws.Range("D91").Value = 85
$ws.Range("D92").Value = 96
$ws.Range("D93").Value = 148
$ws.Range("D94").Value = 180
$ws.Range("D95").Value = 167
$ws.Range("D96").Value = 167
$ws.Range("D97").Value = 135
$ws.Range("D98").Value = 163
$ws.Range("D99").Value = 137
$ws.Range("D100").Value = 146
$ws.Range("D101").Value = 152
$ws.Range("D102").Value = 172
$ws.Range("D103").Value = 145
$ws.Range("D104").Value = 103
$ws.Range("D105").Value = 69
$ws.Range("D106").Value = 32
$ws.Range("D107").Value = 18
$ws.Range("D109").Value = 28
$ws.Range("D110").Value = 30
$ws.Range("D111").Value = 28
$ws.Range("D112").Value = 25
$ws.Range("D113").Value = 22
$ws.Range("D114").Value = 20
$ws.Range("D115").Value = 23
$ws.Range("D116").Value = 26
$ws.Range("D117").Value = 21
$ws.Range("D118").Value = 29
$ws.Range("D119").Value = 19
$ws.Range("D120").Value = 8
$ws.Range("D125").Value = 1
$ws.Range("D127").Value = 1
$ws.Range("D128").Value = 1
$ws.Range("D129").Value = 1
$ws.Range("D131").Value = 1
$ws.Range("D134").Value = 1
$ws.Range("D135").Value = 1
$ws.Range("D140").Value = 7
$ws.Range("D141").Value = 18
$ws.Range("D142").Value = 19
$ws.Range("D143").Value = 24
$ws.Range("D144").Value = 17
$ws.Range("D145").Value = 13
$ws.Range("D146").Value = 17
$ws.Range("D147").Value = 21
$ws.Range("D148").Value = 21
$ws.Range("D149").Value = 37
$ws.Range("D150").Value = 27
$ws.Range("D151").Value = 20
$ws.Range("D152").Value = 11
$ws.Range("D153").Value = 5
$ws.Range("D154").Value = 2
